$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 347.66666
$ws.Range("I4").Value = 347.66666
$ws.Range("K4").Value = 347.66666
$ws.Range("M4").Value = -233.66666

$ws.Range("H5").Value = 456
$ws.Range("I5").Value = 70
$ws.Range("K5").Value = 70
$ws.Range("M5").Value = 45

$ws.Range("H43").Value = 8000
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -7931
$ws.Range("N43").ClearContents()

$ws.Range("H53").Value = 175.81818
$ws.Range("J53").Value = 85.5
$ws.Range("L53").Value = 85.5
$ws.Range("N53").Value = -1359.5

$ws.Range("H76").Value = 5985.2856
$ws.Range("I76").Value = 3849.25
$ws.Range("K76").Value = 3849.25
$ws.Range("M76").Value = -3534.25

$ws.Range("H79").Value = 5985.2856
$ws.Range("I79").Value = 3849.25
$ws.Range("K79").Value = 3849.25
$ws.Range("M79").Value = -2757.25

$ws.Range("H100").Value = 1465.5555
$ws.Range("I100").Value = 1838.8
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 1838.8
$ws.Range("L100").Value = 999
$ws.Range("M100").Value = -1297.8
$ws.Range("N100").Value = -2081

$ws.Range("H111").Value = 1438.3334
$ws.Range("I111").Value = 1280.75
$ws.Range("K111").Value = 3842.25
$ws.Range("M111").Value = -775.25

$ws.Range("H116").Value = 7903
$ws.Range("I116").Value = 9932.611000000001
$ws.Range("K116").Value = 9932.611000000001
$ws.Range("M116").Value = -6490.611000000001

$ws.Range("H132").Value = 2297.0386
$ws.Range("I132").Value = 2129.6667
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6389.000100000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3859.000100000001
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 2393.889
$ws.Range("I137").Value = 1306
$ws.Range("K137").Value = 3918
$ws.Range("M137").Value = -1368

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 21498.5
$ws.Range("J46").Value = 21498.5
$ws.Range("L46").Value = 21498.5
$ws.Range("N46").Value = -22136.5

$ws.Range("H61").Value = 111114650
$ws.Range("I61").Value = 125003416
$ws.Range("K61").Value = 125003416
$ws.Range("M61").Value = -125003204

$ws.Range("H74").Value = 83342260
$ws.Range("I74").Value = 83342260
$ws.Range("K74").Value = 83342260
$ws.Range("M74").Value = -83341386

$ws.Range("H77").Value = 83342260
$ws.Range("I77").Value = 83342260
$ws.Range("K77").Value = 416711300
$ws.Range("M77").Value = -416706932

$ws.Range("H136").Value = 111114650
$ws.Range("I136").Value = 125003416
$ws.Range("K136").Value = 375010248
$ws.Range("M136").Value = -375007698

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4196.8335
$ws.Range("I86").Value = 4340.2
$ws.Range("K86").Value = 4340.2
$ws.Range("M86").Value = -3217.2

$ws.Range("H89").Value = 4196.8335
$ws.Range("I89").Value = 4340.2
$ws.Range("K89").Value = 21701
$ws.Range("M89").Value = -16085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13736
$ws.Range("I31").Value = 8947.9
$ws.Range("J31").Value = 17419.154
$ws.Range("K31").Value = 8947.9
$ws.Range("L31").Value = 17419.154
$ws.Range("M31").Value = -8652.9
$ws.Range("N31").Value = -18009.154

$ws.Range("H34").Value = 13736
$ws.Range("I34").Value = 8947.9
$ws.Range("J34").Value = 17419.154
$ws.Range("K34").Value = 8947.9
$ws.Range("L34").Value = 17419.154
$ws.Range("M34").Value = -8745.9
$ws.Range("N34").Value = -17823.154

$ws.Range("H94").Value = 4500
$ws.Range("J94").Value = 4500
$ws.Range("L94").Value = 4500
$ws.Range("N94").Value = -5402

$ws.Range("H99").Value = 2230.5
$ws.Range("I99").Value = 2013
$ws.Range("K99").Value = 2013
$ws.Range("M99").Value = -515

$ws.Range("H126").Value = 2230.5
$ws.Range("I126").Value = 2013
$ws.Range("K126").Value = 6039
$ws.Range("M126").Value = -3569

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 135.14285
$ws.Range("I2").Value = 120.4
$ws.Range("K2").Value = 722.4000000000001
$ws.Range("M2").Value = -609.4000000000001

$ws.Range("H3").Value = 14843.333
$ws.Range("I3").Value = 5765
$ws.Range("K3").Value = 17295
$ws.Range("M3").Value = -17183

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H107").Value = 2306.9
$ws.Range("I107").Value = 687
$ws.Range("J107").Value = 3386.8333
$ws.Range("K107").Value = 2061
$ws.Range("L107").Value = 10160.4999
$ws.Range("M107").Value = -141
$ws.Range("N107").Value = -14000.4999

$ws.Range("H132").Value = 1529.2222
$ws.Range("J132").Value = 1795.8
$ws.Range("L132").Value = 16162.2
$ws.Range("N132").Value = -21222.2

$ws.Range("H137").Value = 16668622
$ws.Range("I137").Value = 20001948
$ws.Range("K137").Value = 60005844
$ws.Range("M137").Value = -60000744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4311775
$ws.Range("I132").Value = 4311775
$ws.Range("K132").Value = 12935325
$ws.Range("M132").Value = -12932795

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3215.5173
$ws.Range("I40").Value = 3083.3333
$ws.Range("K40").Value = 3083.3333
$ws.Range("M40").Value = -2947.3333

$ws.Range("H93").Value = 5599.6665
$ws.Range("I93").Value = 5599.6665
$ws.Range("K93").Value = 5599.6665
$ws.Range("M93").Value = -4351.6665

$ws.Range("H119").Value = 200420
$ws.Range("J119").Value = 200420
$ws.Range("L119").Value = 200420
$ws.Range("N119").Value = -210096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36144

$ws.Range("H122").Value = 590326.6
$ws.Range("I122").Value = 590326.6
$ws.Range("K122").Value = 1770979.8
$ws.Range("M122").Value = -1768529.8
